$d = $word.ActiveDocument

$replacements = @(
    @("11×29=319", "90×50=4500"),
    @("15×15=225", "60×69=4140"),
    @("77×84=6468", "95×65=6175"),
    @("50×63=3150", "25×87=2175"),
    @("14×30=420", "62×51=3162"),
    @("11×97=1067", "96×61=5856"),
    @("71×78=5538", "25×15=375"),
    @("79×53=4187", "90×98=8820"),
    @("37×88=3256", "73×24=1752"),
    @("93×66=6138", "75×24=1800"),
    @("32×43=1376", "50×73=3650"),
    @("13×38=494", "27×65=1755"),
    @("65×90=5850", "65×28=1820"),
    @("37×68=2516", "11×43=473"),
    @("70×77=5390", "19×68=1292"),
    @("84×43=3612", "67×49=3283"),
    @("98×84=8232", "73×43=3139"),
    @("20×17=340", "86×78=6708"),
    @("94×16=1504", "21×71=1491"),
    @("51×71=3621", "31×56=1736"),
    @("59×30=1770", "33×73=2409"),
    @("40×39=1560", "45×78=3510"),
    @("29×71=2059", "97×91=8827"),
    @("51×37=1887", "47×26=1222"),
    @("82×57=4674", "60×36=2160")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
